# Updates crypto price (D) and 1h-volume-change (E) columns to match the
# latest scrape, per the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold plain numeric-looking strings in the source
# workbook (e.g. "556.52"); Excel's COM Value setter auto-coerces those to
# real numbers, so we force the Text format first and restore the default
# (unstyled) cell format afterwards to avoid introducing style changes.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "59.799.72"
$ws.Range("E2").Value = "  +0.20%  "
Set-TextValue $ws.Range("D3") "2.379.90"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "556.80"
$ws.Range("E5").Value = "  +1.09%  "
Set-TextValue $ws.Range("D6") "133.65"
$ws.Range("E6").Value = "  -2.25%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  +0.34%  "
Set-TextValue $ws.Range("D10") "5.63"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  +1.38%  "
Set-TextValue $ws.Range("D12") "0.343"
$ws.Range("E12").Value = "  -2.83%  "
Set-TextValue $ws.Range("D13") "24.44"
$ws.Range("E13").Value = "  -3.47%  "
Set-TextValue $ws.Range("D14") "2.804.87"
$ws.Range("E14").Value = "  -0.78%  "
Set-TextValue $ws.Range("D15") "59.748.10"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  +0.18%  "
Set-TextValue $ws.Range("D17") "2.383.00"
$ws.Range("E17").Value = "  -0.53%  "
Set-TextValue $ws.Range("D18") "11.12"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +2.21%  "
Set-TextValue $ws.Range("D20") "321.28"
$ws.Range("E20").Value = "  -2.10%  "
Set-TextValue $ws.Range("D21") "6.68"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +0.09%  "
Set-TextValue $ws.Range("D23") "64.11"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +2.11%  "
Set-TextValue $ws.Range("D29") "0.0₃0759"
$ws.Range("E29").Value = "  -1.05%  "
Set-TextValue $ws.Range("D30") "169.84"
$ws.Range("E30").Value = "  +1.02%  "
Set-TextValue $ws.Range("D31") "6.06"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +8.61%  "
$ws.Range("E33").Value = "  -2.40%  "
Set-TextValue $ws.Range("D34") "18.18"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  +0.01%  "
Set-TextValue $ws.Range("D36") "1.33"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E37").Value = "  +0.07%  "
Set-TextValue $ws.Range("D38") "4.13"
$ws.Range("E38").Value = "  -1.01%  "
Set-TextValue $ws.Range("D39") "319.33"
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  -0.85%  "
Set-TextValue $ws.Range("D41") "38.63"
$ws.Range("E41").Value = "  -1.80%  "
Set-TextValue $ws.Range("D42") "145.31"
$ws.Range("E42").Value = "  +4.58%  "
$ws.Range("E43").Value = "  -3.76%  "
Set-TextValue $ws.Range("D44") "0.0970"
$ws.Range("E44").Value = "  +0.19%  "
Set-TextValue $ws.Range("D45") "19.67"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  -2.18%  "
Set-TextValue $ws.Range("D49") "11.07"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -1.66%  "
